$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "60.279.06"
$ws.Range("E2").Value = "  +3.72%  "

$ws.Range("D3").Value = "2.332.19"
$ws.Range("E3").Value = "  +2.02%  "

$ws.Range("E4").Value = "  -0.04%  "

Set-TextValue $ws.Range("D5") "547.95"
$ws.Range("E5").Value = "  +2.63%  "

Set-TextValue $ws.Range("D6") "131.58"
$ws.Range("E6").Value = "  +0.66%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  -0.62%  "

$ws.Range("D9").Value = "2.329.52"
$ws.Range("E9").Value = "  +2.00%  "

$ws.Range("E10").Value = "  +1.44%  "

Set-TextValue $ws.Range("D11") "5.53"
$ws.Range("E11").Value = "  +0.86%  "

$ws.Range("E12").Value = "  +0.40%  "

Set-TextValue $ws.Range("D13") "0.336"
$ws.Range("E13").Value = "  +1.52%  "

Set-TextValue $ws.Range("D14") "23.78"
$ws.Range("E14").Value = "  +1.64%  "

$ws.Range("D15").Value = "2.746.85"
$ws.Range("E15").Value = "  +1.92%  "

$ws.Range("D16").Value = "60.235.58"
$ws.Range("E16").Value = "  +3.75%  "

$ws.Range("E17").Value = "  +1.07%  "

$ws.Range("D18").Value = "2.330.50"
$ws.Range("E18").Value = "  +2.53%  "

Set-TextValue $ws.Range("D19") "10.63"
$ws.Range("E19").Value = "  +1.41%  "

$ws.Range("E20").Value = "  -0.53%  "

Set-TextValue $ws.Range("D21") "314.91"
$ws.Range("E21").Value = "  +0.54%  "

$ws.Range("E22").Value = "  +4.11%  "

$ws.Range("E23").Value = "  -0.32%  "

Set-TextValue $ws.Range("D24") "64.05"
$ws.Range("E24").Value = "  +2.22%  "

Set-TextValue $ws.Range("D25") "0.171"
$ws.Range("E25").Value = "  +1.78%  "

$ws.Range("E26").Value = "  +0.34%  "

Set-TextValue $ws.Range("D27") "7.86"
$ws.Range("E27").Value = "  -1.19%  "

Set-TextValue $ws.Range("D28") "1.35"
$ws.Range("E28").Value = "  +7.49%  "

Set-TextValue $ws.Range("D29") "1.22"
$ws.Range("E29").Value = "  +16.38%  "

Set-TextValue $ws.Range("D30") "174.32"
$ws.Range("E30").Value = "  +1.98%  "

$ws.Range("E31").Value = "  +3.27%  "

$ws.Range("D32").Value = "0.0₃0732"
$ws.Range("E32").Value = "  +1.58%  "

Set-TextValue $ws.Range("D33") "5.96"
$ws.Range("E33").Value = "  +3.65%  "

$ws.Range("E34").Value = "  +11.68%  "

$ws.Range("E35").Value = "  +0.71%  "

Set-TextValue $ws.Range("D37") "17.91"
$ws.Range("E37").Value = "  +0.72%  "

$ws.Range("E38").Value = "  -0.12%  "

$ws.Range("E39").Value = "  +5.21%  "

Set-TextValue $ws.Range("D40") "334.38"
$ws.Range("E40").Value = "  +16.55%  "

Set-TextValue $ws.Range("D41") "38.02"
$ws.Range("E41").Value = "  -0.87%  "

Set-TextValue $ws.Range("D42") "1.53"
$ws.Range("E42").Value = "  +2.77%  "

Set-TextValue $ws.Range("D43") "139.96"
$ws.Range("E43").Value = "  -0.44%  "

Set-TextValue $ws.Range("D44") "3.48"
$ws.Range("E44").Value = "  +1.53%  "

Set-TextValue $ws.Range("D45") "0.0944"
$ws.Range("E45").Value = "  -0.75%  "

Set-TextValue $ws.Range("D46") "19.34"
$ws.Range("E46").Value = "  +8.05%  "

Set-TextValue $ws.Range("D47") "0.0497"
$ws.Range("E47").Value = "  +0.97%  "

Set-TextValue $ws.Range("D48") "0.560"
$ws.Range("E48").Value = "  +1.55%  "

$ws.Range("D49").Value = "0.0₆0228"
$ws.Range("E49").Value = "  +22.16%  "

Set-TextValue $ws.Range("D50") "0.0214"
$ws.Range("E50").Value = "  +2.01%  "

$ws.Range("E51").Value = "  +0.66%  "

